$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "set" row (row 2) and the "banana" row (row 4), keeping the
# "apple" row (row 3, which shifts up to row 2). The apple row already has
# an empty-string shared-string cell in its Antonyms column, which we want
# to preserve as-is in the final "words" row.
$ws.Rows("4:4").Delete()
$ws.Rows("2:2").Delete()

$wordsJson = @'
{"definition":"The smallest unit of language that has a particular meaning and can be expressed by itself; the smallest discrete, meaningful unit of language. (contrast morpheme.)","synonyms":[],"antonyms":[]},{"definition":"Something like such a unit of language:","synonyms":[],"antonyms":[]},{"definition":"The fact or act of speaking, as opposed to taking action. .","synonyms":[],"antonyms":[]},{"definition":"Something that someone said; a comment, utterance; speech.","synonyms":[],"antonyms":[]},{"definition":"A watchword or rallying cry, a verbal signal (even when consisting of multiple words).","synonyms":[],"antonyms":[],"example":"mum's the word"},{"definition":"A proverb or motto.","synonyms":[],"antonyms":[]},{"definition":"News; tidings (used without an article).","synonyms":[],"antonyms":[],"example":"Have you had any word from John yet?"},{"definition":"An order; a request or instruction; an expression of will.","synonyms":[],"antonyms":[],"example":"Don't fire till I give the word"},{"definition":"A promise; an oath or guarantee.","synonyms":["promise"],"antonyms":[],"example":"I give you my word that I will be there on time."},{"definition":"A brief discussion or conversation.","synonyms":[],"antonyms":[],"example":"Can I have a word with you?"},{"definition":"(in the plural) See words.","synonyms":[],"antonyms":[],"example":"There had been words between him and the secretary about the outcome of the meeting."},{"definition":"(sometimes Word) Communication from God; the message of the Christian gospel; the Bible, Scripture.","synonyms":["Bible","word of God"],"antonyms":[],"example":"Her parents had lived in Botswana, spreading the word among the tribespeople."},{"definition":"(sometimes Word) Logos, Christ.","synonyms":["God","Logos"],"antonyms":[]}
'@

# Overwrite the remaining data row (now row 2) with the "words" entry. The
# Antonyms cell (F2) is intentionally left untouched so it keeps its
# existing empty shared-string value instead of being cleared.
$ws.Range("A2").Value = "words"
$ws.Range("B2").Value = "noun"
$ws.Range("C2").Value = "/wɜːdz/"
$ws.Range("D2").Value = $wordsJson
$ws.Range("E2").Value = "Bible,word of God,God,Logos,promise,vocable"
